$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '58.249.08'
$dCell.ClearFormats()
$ws.Range("E2").Value = '  +0.68%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '3.152.23'
$dCell.ClearFormats()
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.01%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '534.93'
$dCell.ClearFormats()
$ws.Range("E5").Value = '  +1.12%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '138.84'
$dCell.ClearFormats()
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.03%  '

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '3.150.93'
$dCell.ClearFormats()
$ws.Range("E8").Value = '  +0.59%  '

$ws.Range("E9").Value = '  +4.42%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '7.32'
$dCell.ClearFormats()
$ws.Range("E10").Value = '  +1.91%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.108'
$dCell.ClearFormats()
$ws.Range("E11").Value = '  -0.11%  '

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '0.413'
$dCell.ClearFormats()
$ws.Range("E12").Value = '  +4.12%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '3.669.53'
$dCell.ClearFormats()
$ws.Range("E13").Value = '  -0.03%  '

$ws.Range("E14").Value = '  +1.44%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '25.82'
$dCell.ClearFormats()
$ws.Range("E15").Value = '  +1.15%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0000165'
$dCell.ClearFormats()
$ws.Range("E16").Value = '  +0.26%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '58.305.12'
$dCell.ClearFormats()
$ws.Range("E17").Value = '  +0.58%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '3.144.05'
$dCell.ClearFormats()
$ws.Range("E18").Value = '  +0.59%  '

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '6.04'
$dCell.ClearFormats()
$ws.Range("E19").Value = '  +0.06%  '

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '12.74'
$dCell.ClearFormats()
$ws.Range("E20").Value = '  -0.72%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '8.19'
$dCell.ClearFormats()
$ws.Range("E21").Value = '  +2.88%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '361.26'
$dCell.ClearFormats()
$ws.Range("E22").Value = '  +2.11%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '0.997'
$dCell.ClearFormats()
$ws.Range("E23").Value = '  -0.43%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '69.23'
$dCell.ClearFormats()
$ws.Range("E24").Value = '  +0.95%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '0.507'
$dCell.ClearFormats()
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  -1.38%  '

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '0.998'
$dCell.ClearFormats()
$ws.Range("E27").Value = '  -0.32%  '

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₃0883'
$dCell.ClearFormats()
$ws.Range("E28").Value = '  -3.92%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '7.35'
$dCell.ClearFormats()
$ws.Range("E29").Value = '  -1.83%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '6.19'
$dCell.ClearFormats()
$ws.Range("E30").Value = '  -0.37%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '1.88'
$dCell.ClearFormats()
$ws.Range("E31").Value = '  -0.35%  '

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '21.54'
$dCell.ClearFormats()
$ws.Range("E32").Value = '  +1.82%  '

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '5.02'
$dCell.ClearFormats()
$ws.Range("E33").Value = '  +0.78%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.15'
$dCell.ClearFormats()
$ws.Range("E34").Value = '  -3.38%  '

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '158.77'
$dCell.ClearFormats()
$ws.Range("E35").Value = '  +0.30%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '6.10'
$dCell.ClearFormats()
$ws.Range("E36").Value = '  -1.24%  '

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '25.93'
$dCell.ClearFormats()
$ws.Range("E37").Value = '  -1.93%  '

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '1.28'
$dCell.ClearFormats()
$ws.Range("E38").Value = '  -0.26%  '

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '1.70'
$dCell.ClearFormats()
$ws.Range("E39").Value = '  +4.91%  '

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0671'
$dCell.ClearFormats()
$ws.Range("E40").Value = '  +0.30%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '2.518.72'
$dCell.ClearFormats()
$ws.Range("E41").Value = '  +8.30%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '0.703'
$dCell.ClearFormats()
$ws.Range("E42").Value = '  -0.21%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '4.02'
$dCell.ClearFormats()
$ws.Range("E43").Value = '  -3.53%  '

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '37.48'
$dCell.ClearFormats()
$ws.Range("E44").Value = '  +2.34%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '3.185.97'
$dCell.ClearFormats()
$ws.Range("E45").Value = '  +0.46%  '

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0270'
$dCell.ClearFormats()
$ws.Range("E46").Value = '  -0.52%  '

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '0.999'
$dCell.ClearFormats()
$ws.Range("E47").Value = '  -0.03%  '

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '0.987'
$dCell.ClearFormats()
$ws.Range("E48").Value = '  +1.51%  '

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '6.06'
$dCell.ClearFormats()
$ws.Range("E49").Value = '  +0.43%  '

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '19.87'
$dCell.ClearFormats()
$ws.Range("E50").Value = '  -3.20%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '0.740'
$dCell.ClearFormats()
$ws.Range("E51").Value = '  -3.88%  '
